$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 9849.31
$ws.Cells.Item(6, 2).Value = 9949.7999999999993
$ws.Cells.Item(6, 3).Value = 286
$ws.Cells.Item(6, 4).Value = 283.11
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(6, 6).Value = -1.01
$ws.Cells.Item(6, 7).Value = 42612.675069444442
$ws.Cells.Item(6, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(6, 8).Value = $false
